$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data blocks for rows 2-5 (date 44908) and rows 6-9 (date 44890)
# Columns affected: D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen), S (Precio $/Kg)

$cols = @("D","M","N","O","P","R","S")
$pairs = @(@(2,6), @(3,7), @(4,8), @(5,9))

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}
